$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    3 = @(0.6753301551942219, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 2.020003488691574)
    4 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    5 = @(0.04763786555579896, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 0.7443468554461139)
    6 = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    7 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    8 = @(0.127881588408715, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1.094976487407548)
    9 = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 8.660232485948974, 9.906394106792828)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
